$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
$ws.Columns.Item(9).ColumnWidth = 2.140625
$ws.Columns.Item(11).ColumnWidth = 5.7109375
$ws.Columns.Item(12).ColumnWidth = 5.7109375
$ws.Columns.Item(14).ColumnWidth = 5.7109375

# Update row 1 values
$ws.Range("B1").Value = 4
$ws.Range("C1").Value = 27
$ws.Range("D1").Value = 11
$ws.Range("E1").Value = 32
$ws.Range("F1").Value = 23
$ws.Range("G1").Value = 20
$ws.Range("H1").Value = 21
$ws.Range("I1").Value = 2
$ws.Range("J1").Value = 15
$ws.Range("K1").Value = 0.045
$ws.Range("L1").Value = 0.019
$ws.Range("M1").Value = 0.053
$ws.Range("N1").Value = 0.097
